# Auto-generated script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.076.18"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "1.897.66"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5025"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  -1.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09190"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.128"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.81"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.392"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").Value = "1.898.69"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("E15").Value = "  -3.47%  "
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.36"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06659"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.92"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.210"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("D23").Value = "28.124.47"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.46"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.586"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.50%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.112.63"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.90"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.19"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.44"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.76%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.090"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1059"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.61%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.607"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.617"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.587"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06602"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.73%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02406"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2208"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.227"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.73%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.274"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.46%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6499"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.977"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.38"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6099"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.32"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.301"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.684"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.001"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.98"
$ws.Range("D50").ClearFormats()
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.183"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.22%  "
